# Update "Forecast Comparison" sheet: insert a Week_Start_Date column,
# normalize the Week labels, refresh MyForecast numbers and mark
# is_holiday_week as a boolean. Also refresh the dependent Summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column before the current column B (ASIN), pushing the
#    old B..I columns to C..J.
$ws.Columns.Item(2).Insert()

# 2. New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# 3. Week start date (one per forecast week) + normalized Week label +
#    refreshed MyForecast value, written per data row (rows 2-17).
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)
$weekLabels = @(
    "W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8",
    "W9", "W10", "W11", "W12", "W13", "W14", "W15", "W16"
)
$myForecast = @(29, 37, 68, 57, 61, 69, 71, 67, 68, 64, 66, 62, 65, 61, 64, 60)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Week label without the leading zero (W01 -> W1, etc.).
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]

    # Week_Start_Date must stay plain text, not get auto-converted to an
    # Excel date serial number.
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]

    # Refreshed MyForecast figure (now in column D after the insert).
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]

    # is_holiday_week (column J after the insert) becomes a boolean.
    $ws.Cells.Item($row, 10).Value = $false
}

# 4. Refresh the dependent totals on the Summary sheet. The Value column on
#    this sheet is stored as plain text (e.g. "3629 units"), so keep these
#    updated numbers as text too instead of letting Excel turn them into
#    numeric cells.
$summary = $wb.Worksheets.Item("Summary")

$summaryUpdates = @{ "B9" = "969"; "B10" = "459"; "B12" = "71" }
foreach ($addr in $summaryUpdates.Keys) {
    $cell = $summary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
}
